$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A48").Value = "2025-04-29 05:33:17"
$ws.Range("B48").Value = 138
